$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date-style formatting (style index used by column A, rows 2-343)
# from the last existing data row (343) down onto the new rows (344-357),
# then fill in the actual data values for 10-23 August 2021.
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)

$ws.Range("A344").Value = 44418
$ws.Range("B344").Value = 2
$ws.Range("C344").Value = 18
$ws.Range("D344").Value = 606.8779501011463

$ws.Range("A345").Value = 44419
$ws.Range("B345").Value = 1
$ws.Range("C345").Value = 18
$ws.Range("D345").Value = 606.8779501011463

$ws.Range("A346").Value = 44420
$ws.Range("B346").Value = 2
$ws.Range("C346").Value = 19
$ws.Range("D346").Value = 640.5933917734322

$ws.Range("A347").Value = 44421
$ws.Range("B347").Value = 2
$ws.Range("C347").Value = 19
$ws.Range("D347").Value = 640.5933917734322

$ws.Range("A348").Value = 44422
$ws.Range("B348").Value = 0
$ws.Range("C348").Value = 15
$ws.Range("D348").Value = 505.7316250842886

$ws.Range("A349").Value = 44423
$ws.Range("B349").Value = 0
$ws.Range("C349").Value = 8
$ws.Range("D349").Value = 269.7235333782872

$ws.Range("A350").Value = 44424
$ws.Range("B350").Value = 2
$ws.Range("C350").Value = 9
$ws.Range("D350").Value = 303.4389750505732

$ws.Range("A351").Value = 44425
$ws.Range("B351").Value = 3
$ws.Range("C351").Value = 10
$ws.Range("D351").Value = 337.1544167228591

$ws.Range("A352").Value = 44426
$ws.Range("B352").Value = 2
$ws.Range("C352").Value = 11
$ws.Range("D352").Value = 370.8698583951449

$ws.Range("A353").Value = 44427
$ws.Range("B353").Value = 3
$ws.Range("C353").Value = 12
$ws.Range("D353").Value = 404.5853000674309

$ws.Range("A354").Value = 44428
$ws.Range("B354").Value = 1
$ws.Range("C354").Value = 11
$ws.Range("D354").Value = 370.8698583951449

$ws.Range("A355").Value = 44429
$ws.Range("B355").Value = 0
$ws.Range("C355").Value = 11
$ws.Range("D355").Value = 370.8698583951449

$ws.Range("A356").Value = 44430
$ws.Range("B356").Value = 0
$ws.Range("C356").Value = 11
$ws.Range("D356").Value = 370.8698583951449

$ws.Range("A357").Value = 44431
$ws.Range("B357").Value = 2
$ws.Range("C357").Value = 11
$ws.Range("D357").Value = 370.8698583951449

